# Insert two blank rows above the existing table. Excel's row insert copies
# the formatting of the row that was previously row 1 (the bold/bordered
# header style) down onto what becomes row 3, and leaves the two newly
# inserted rows (1 and 2) unformatted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:2").Insert()

# Row 3 now holds the original header text but still carries the old
# header's bold/border/center style. Copy that style up onto the new row 1
# (which will hold the numeric column-index row) before stripping it from
# row 3, so row 3 ends up with default formatting like a normal data row.
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A1:N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A3:N3").ClearFormats()

# Row 3's K/M/N cells (part number / thread size / material columns in the
# old header) are blanked out in the new layout.
$ws.Range("M3:N3").ClearContents()

# New row 1: sequential numeric column indexes 0-13.
$colIndexValues = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13)
for ($i = 0; $i -lt $colIndexValues.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $colIndexValues[$i]
}

# New row 2: only column E gets a label, the rest stays blank.
$ws.Range("E2").Value = "Washer"
